# This workbook ("Hortaliza, Vega Modelo de Temuco - Puerro") tracks daily
# price observations for Puerro (leek) at the Vega Modelo de Temuco market.
# A new observation was inserted as row 202 (pushing the existing rows 202-253
# down to 203-254), growing the used range from A1:R253 to A1:R254.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 202; this shifts rows 202:253 down to 203:254
# and extends the sheet dimension automatically.
$ws.Rows.Item(202).Insert()

# Populate the newly inserted row 202 with the new observation's data.
$ws.Cells.Item(202, 1).Value  = 10
$ws.Cells.Item(202, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(202, 3).Value  = "La Araucanía"
$ws.Cells.Item(202, 4).Value  = 44889
$ws.Cells.Item(202, 5).Value  = 9
$ws.Cells.Item(202, 6).Value  = 100112005
$ws.Cells.Item(202, 7).Value  = "Puerro"
$ws.Cells.Item(202, 8).Value  = "Azul de Maquehue"
$ws.Cells.Item(202, 9).Value  = "Primera"
$ws.Cells.Item(202, 10).Value = 65
$ws.Cells.Item(202, 11).Value = 18000
$ws.Cells.Item(202, 12).Value = 18000
$ws.Cells.Item(202, 13).Value = 18000
$ws.Cells.Item(202, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(202, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(202, 16).Value = 1500
$ws.Cells.Item(202, 17).Value = 12
$ws.Cells.Item(202, 18).Value = "Hortaliza"
